# Apply the benchmark-stat corrections described in the commit.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

# Simple value replacements (single run, no restructuring needed).
Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "459"
Set-CellText $t 6  "0.01910"
Set-CellText $t 7  "0.00164"
Set-CellText $t 8  "0.00063"
Set-CellText $t 9  "0.00804"
Set-CellText $t 10 "0.01009"
Set-CellText $t 11 "0.01286"
Set-CellText $t 12 "0.07936"

# Collapse the three multi-run / tab-delimited rows down to a single value.
Set-CellText $t 44 "99.93"
Set-CellText $t 45 "0.08"
Set-CellText $t 46 "119"
